$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.799369692802429
$ws.Range("B1").Value = 2.312338590621948
$ws.Range("C1").Value = 1.99890124797821
$ws.Range("D1").Value = 1.663317203521729
$ws.Range("E1").Value = 1.575534820556641
